# Fruta / hortaliza, semanal
# Insert two new weekly rows (483 and 484) above the existing data block,
# shifting the former rows 483:519 down to 485:521.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 483 (pushes 483-519 down to 485-521,
# and extends the sheet dimension to A1:T521 automatically).
$ws.Range("A483:A484").EntireRow.Insert()

# --- New row 483 ---
$ws.Range("A483").Value = 4
$ws.Range("B483").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C483").Value = "Los Lagos"
$ws.Range("D483").Value = 45265
$ws.Range("E483").Value = 10
$ws.Range("F483").Value = "Fruta"
$ws.Range("G483").Value = 100101
$ws.Range("H483").Value = "Berries"
$ws.Range("I483").Value = 100112025
$ws.Range("J483").Value = "Frutilla"
$ws.Range("K483").Value = "Sin especificar"
$ws.Range("L483").Value = "Primera"
$ws.Range("M483").Value = 300
$ws.Range("N483").Value = 15000
$ws.Range("O483").Value = 15000
$ws.Range("P483").Value = 15000
$ws.Range("Q483").Value = "$/bandeja 7 kilos"
$ws.Range("R483").Value = "Provincia de Melipilla"
$ws.Range("S483").Value = 2143
$ws.Range("T483").Value = 7

# --- New row 484 ---
$ws.Range("A484").Value = 4
$ws.Range("B484").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C484").Value = "Los Lagos"
$ws.Range("D484").Value = 45265
$ws.Range("E484").Value = 10
$ws.Range("F484").Value = "Fruta"
$ws.Range("G484").Value = 100101
$ws.Range("H484").Value = "Berries"
$ws.Range("I484").Value = 100112025
$ws.Range("J484").Value = "Frutilla"
$ws.Range("K484").Value = "Sin especificar"
$ws.Range("L484").Value = "Primera"
$ws.Range("M484").Value = 300
$ws.Range("N484").Value = 14000
$ws.Range("O484").Value = 14000
$ws.Range("P484").Value = 14000
$ws.Range("Q484").Value = "$/caja 7 kilos"
$ws.Range("R484").Value = "Región de La Araucanía"
$ws.Range("S484").Value = 2000
$ws.Range("T484").Value = 7
